$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.20127533333333
$ws.Range("H2").Value = 33.603826
$ws.Range("I2").Value = 0.1186573945858706
$ws.Range("J2").Value = 0.1186573945858706
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 192.8285726666667
$ws.Range("N2").Value = 578.485718
$ws.Range("O2").Value = 0.7801188850698786
$ws.Range("P2").Value = 0.7801188850698786
$ws.Range("Q2").Value = 2159.925934573007
$ws.Range("R2").Value = 19439.33341115707
$ws.Range("S2").Value = 0.09256687436962605
$ws.Range("T2").Value = 0.09256687436962605

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.20127533333333
$ws.Range("H3").Value = 33.603826
$ws.Range("I3").Value = 0.1186573945858706
$ws.Range("J3").Value = 0.1186573945858706
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4209206666666667
$ws.Range("N3").Value = 1.262762
$ws.Range("O3").Value = 0.001702901995496819
$ws.Range("P3").Value = 0.001702901995496819
$ws.Range("Q3").Value = 4.714848280823556
$ws.Range("R3").Value = 42.433634527412
$ws.Range("S3").Value = 0.0002020619140207325
$ws.Range("T3").Value = 0.0002020619140207325

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.20127533333333
$ws.Range("H4").Value = 33.603826
$ws.Range("I4").Value = 0.1186573945858706
$ws.Range("J4").Value = 0.1186573945858706
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 45.70525533333333
$ws.Range("N4").Value = 137.115766
$ws.Range("O4").Value = 0.184907933193646
$ws.Range("P4").Value = 0.184907933193646
$ws.Range("Q4").Value = 511.9571491689684
$ws.Range("R4").Value = 4607.614342520716
$ws.Range("S4").Value = 0.02194069359101626
$ws.Range("T4").Value = 0.02194069359101626

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.20127533333333
$ws.Range("H5").Value = 33.603826
$ws.Range("I5").Value = 0.1186573945858706
$ws.Range("J5").Value = 0.1186573945858706
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.223696
$ws.Range("N5").Value = 24.671088
$ws.Range("O5").Value = 0.0332702797409786
$ws.Range("P5").Value = 0.0332702797409786
$ws.Range("Q5").Value = 92.115883153632
$ws.Range("R5").Value = 829.042948382688
$ws.Range("S5").Value = 0.003947764711207596
$ws.Range("T5").Value = 0.003947764711207595

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("H6").Value = 143.833961
$ws.Range("I6").Value = 0.5078874966566524
$ws.Range("J6").Value = 0.5078874966566524
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 192.8285726666667
$ws.Range("N6").Value = 578.485718
$ws.Range("O6").Value = 0.7801188850698786
$ws.Range("P6").Value = 0.7801188850698786
$ws.Range("Q6").Value = 9245.099133541002
$ws.Range("R6").Value = 83205.89220186901
$ws.Range("S6").Value = 0.3962126276327194
$ws.Range("T6").Value = 0.3962126276327194

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("H7").Value = 143.833961
$ws.Range("I7").Value = 0.5078874966566524
$ws.Range("J7").Value = 0.5078874966566524
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4209206666666667
$ws.Range("N7").Value = 1.262762
$ws.Range("O7").Value = 0.001702901995496819
$ws.Range("P7").Value = 0.001702901995496819
$ws.Range("Q7").Value = 20.18089558447578
$ws.Range("R7").Value = 181.628060260282
$ws.Range("S7").Value = 0.0008648826315444974
$ws.Range("T7").Value = 0.0008648826315444973

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("H8").Value = 143.833961
$ws.Range("I8").Value = 0.5078874966566524
$ws.Range("J8").Value = 0.5078874966566524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.70525533333333
$ws.Range("N8").Value = 137.115766
$ws.Range("O8").Value = 0.184907933193646
$ws.Range("P8").Value = 0.184907933193646
$ws.Range("Q8").Value = 2191.322637703237
$ws.Range("R8").Value = 19721.90373932913
$ws.Range("S8").Value = 0.09391242730167641
$ws.Range("T8").Value = 0.09391242730167641

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("H9").Value = 143.833961
$ws.Range("I9").Value = 0.5078874966566524
$ws.Range("J9").Value = 0.5078874966566524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.223696
$ws.Range("N9").Value = 24.671088
$ws.Range("O9").Value = 0.0332702797409786
$ws.Range("P9").Value = 0.0332702797409786
$ws.Range("Q9").Value = 394.2822565799521
$ws.Range("R9").Value = 3548.540309219568
$ws.Range("S9").Value = 0.01689755909071216
$ws.Range("T9").Value = 0.01689755909071216

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.59984766666667
$ws.Range("H10").Value = 76.799543
$ws.Range("I10").Value = 0.2711844085184091
$ws.Range("J10").Value = 0.2711844085184091
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 192.8285726666667
$ws.Range("N10").Value = 578.485718
$ws.Range("O10").Value = 0.7801188850698786
$ws.Range("P10").Value = 0.7801188850698786
$ws.Range("Q10").Value = 4936.382086047431
$ws.Range("R10").Value = 44427.43877442688
$ws.Range("S10").Value = 0.2115560784217158
$ws.Range("T10").Value = 0.2115560784217158

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 25.59984766666667
$ws.Range("H11").Value = 76.799543
$ws.Range("I11").Value = 0.2711844085184091
$ws.Range("J11").Value = 0.2711844085184091
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.4209206666666667
$ws.Range("N11").Value = 1.262762
$ws.Range("O11").Value = 0.001702901995496819
$ws.Range("P11").Value = 0.001702901995496819
$ws.Range("Q11").Value = 10.77550494641845
$ws.Range("R11").Value = 96.97954451776602
$ws.Range("S11").Value = 0.0004618004704136235
$ws.Range("T11").Value = 0.0004618004704136234

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 25.59984766666667
$ws.Range("H12").Value = 76.799543
$ws.Range("I12").Value = 0.2711844085184091
$ws.Range("J12").Value = 0.2711844085184091
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 45.70525533333333
$ws.Range("N12").Value = 137.115766
$ws.Range("O12").Value = 0.184907933193646
$ws.Range("P12").Value = 0.184907933193646
$ws.Range("Q12").Value = 1170.047574099437
$ws.Range("R12").Value = 10530.42816689494
$ws.Range("S12").Value = 0.05014414849348042
$ws.Range("T12").Value = 0.05014414849348042

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 25.59984766666667
$ws.Range("H13").Value = 76.799543
$ws.Range("I13").Value = 0.2711844085184091
$ws.Range("J13").Value = 0.2711844085184091
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.223696
$ws.Range("N13").Value = 24.671088
$ws.Range("O13").Value = 0.0332702797409786
$ws.Range("P13").Value = 0.0332702797409786
$ws.Range("Q13").Value = 210.525364856976
$ws.Range("R13").Value = 1894.728283712784
$ws.Range("S13").Value = 0.009022381132799294
$ws.Range("T13").Value = 0.009022381132799292

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.654369000000001
$ws.Range("H14").Value = 28.963107
$ws.Range("I14").Value = 0.1022707002390678
$ws.Range("J14").Value = 0.1022707002390678
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 192.8285726666667
$ws.Range("N14").Value = 578.485718
$ws.Range("O14").Value = 0.7801188850698786
$ws.Range("P14").Value = 0.7801188850698786
$ws.Range("Q14").Value = 1861.638194267314
$ws.Range("R14").Value = 16754.74374840583
$ws.Range("S14").Value = 0.07978330464581733
$ws.Range("T14").Value = 0.07978330464581732

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.654369000000001
$ws.Range("H15").Value = 28.963107
$ws.Range("I15").Value = 0.1022707002390678
$ws.Range("J15").Value = 0.1022707002390678
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.4209206666666667
$ws.Range("N15").Value = 1.262762
$ws.Range("O15").Value = 0.001702901995496819
$ws.Range("P15").Value = 0.001702901995496819
$ws.Range("Q15").Value = 4.063723435726001
$ws.Range("R15").Value = 36.57351092153401
$ws.Range("S15").Value = 0.0001741569795179655
$ws.Range("T15").Value = 0.0001741569795179655

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.654369000000001
$ws.Range("H16").Value = 28.963107
$ws.Range("I16").Value = 0.1022707002390678
$ws.Range("J16").Value = 0.1022707002390678
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 45.70525533333333
$ws.Range("N16").Value = 137.115766
$ws.Range("O16").Value = 0.184907933193646
$ws.Range("P16").Value = 0.184907933193646
$ws.Range("Q16").Value = 441.255400227218
$ws.Range("R16").Value = 3971.298602044963
$ws.Range("S16").Value = 0.01891066380747295
$ws.Range("T16").Value = 0.01891066380747294

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.654369000000001
$ws.Range("H17").Value = 28.963107
$ws.Range("I17").Value = 0.1022707002390678
$ws.Range("J17").Value = 0.1022707002390678
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.223696
$ws.Range("N17").Value = 24.671088
$ws.Range("O17").Value = 0.0332702797409786
$ws.Range("P17").Value = 0.0332702797409786
$ws.Range("Q17").Value = 79.39459572782401
$ws.Range("R17").Value = 714.551361550416
$ws.Range("S17").Value = 0.003402574806259553
$ws.Range("T17").Value = 0.003402574806259552
